$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Sheet 2")

$ws.Range("C14").Value = "Compare Val 1"
$ws.Range("D14").Value = 230

$ws.Range("C15").Value = "Compare Val 2"
$ws.Range("D15").Value = 230

$ws.Range("D15").Select()
